# document_loader/src/test.xlsx — add a second sheet ("Sheet2") that mirrors
# Sheet1's layout but with its own header labels and an incrementing id of 2
# instead of 1, then leave Sheet2 selected/active as the new tab.

$wb = $excel.ActiveWorkbook

# --- Sheet1: existing worksheet -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet2: new worksheet, inserted right after Sheet1 -------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row. Written B, A, C so the new shared-string table entries land in
# the same order as the target file (Test4, Test3, Time-2).
$ws2.Cells.Item(1, 2).Value = "Test4"
$ws2.Cells.Item(1, 1).Value = "Test3"
$ws2.Cells.Item(1, 3).Value = "Time-2"

# Data rows 2..22: column A is a constant 2, column B counts down from 100,
# column C counts up from 0 (same shape as Sheet1, whose column A is 1).
for ($r = 2; $r -le 22; $r++) {
    $ws2.Cells.Item($r, 1).Value = 2
    $ws2.Cells.Item($r, 2).Value = 102 - $r
    $ws2.Cells.Item($r, 3).Value = $r - 2
}

# --- View state -------------------------------------------------------------
# Sheet1: scroll down so row 6 is at the top, select the whole used range.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 6
$ws1.Range("A1:C22").Select()

# Sheet2: becomes the active/visible tab, scrolled the same way, with B2:B22
# equivalent selection shifted onto column A (A2:A22, active cell A2).
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 6
$ws2.Range("A2:A22").Select()
